$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new inventory row (row 32) matching the pattern of existing rows
$ws.Range("A32").Value = "UAYYDQ"
$ws.Range("B32").Value = "DMD para proyector 8060-6039B"
$ws.Range("C32").Value = "Benq MP515 MP515ST NEC NP115 OPTOMA ES526"
$ws.Range("D32").Value = 200000
$ws.Range("E32").Value = 400000
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 1
$ws.Range("H32").Formula = "=(E32-D32)*G32"
$ws.Range("I32").Formula = "=D32*F32"
$ws.Range("J32").Value = 200000
